# Auto-generated edit script applying commit 'Correcion a Diebold Mariano y revision de Cap1'
# Updates cell values in sheets: Matriz_Resultados, P_valores, Estadisticos_DM, Resumen
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Matriz_Resultados")
$ws.Range("C2").Value = [double]"0"
$ws.Range("I2").Value = [double]"0"
$ws.Range("B3").Value = [double]"0"
$ws.Range("E3").Value = [double]"0"
$ws.Range("F3").Value = [double]"0"
$ws.Range("I3").Value = [double]"0"
$ws.Range("E4").Value = [double]"0"
$ws.Range("C5").Value = [double]"0"
$ws.Range("D5").Value = [double]"0"
$ws.Range("G5").Value = [double]"0"
$ws.Range("H5").Value = [double]"0"
$ws.Range("I5").Value = [double]"0"
$ws.Range("C6").Value = [double]"0"
$ws.Range("I6").Value = [double]"0"
$ws.Range("E7").Value = [double]"0"
$ws.Range("E8").Value = [double]"0"
$ws.Range("B9").Value = [double]"0"
$ws.Range("C9").Value = [double]"0"
$ws.Range("E9").Value = [double]"0"
$ws.Range("F9").Value = [double]"0"

$ws = $wb.Worksheets.Item("P_valores")
$ws.Range("C2").Value = [double]"0.003939399324226045"
$ws.Range("D2").Value = [double]"0.04108229454594259"
$ws.Range("E2").Value = [double]"0.2502371920658693"
$ws.Range("F2").Value = [double]"0.6387018093091172"
$ws.Range("G2").Value = [double]"0.03847178559399977"
$ws.Range("H2").Value = [double]"0.02324467905162675"
$ws.Range("I2").Value = [double]"0.004331495970726795"
$ws.Range("J2").Value = [double]"0.504032656816424"
$ws.Range("B3").Value = [double]"0.003939399324226045"
$ws.Range("D3").Value = [double]"6.137391683758153E-08"
$ws.Range("E3").Value = [double]"0.003229026081643438"
$ws.Range("F3").Value = [double]"0.003849857549612201"
$ws.Range("G3").Value = [double]"7.060412010595485E-07"
$ws.Range("H3").Value = [double]"0.0001782056793002695"
$ws.Range("I3").Value = [double]"0.006419911705607184"
$ws.Range("J3").Value = [double]"4.817164197534751E-06"
$ws.Range("B4").Value = [double]"0.04108229454594259"
$ws.Range("C4").Value = [double]"6.137391683758153E-08"
$ws.Range("E4").Value = [double]"0.02011128779015348"
$ws.Range("F4").Value = [double]"0.03618311219722914"
$ws.Range("G4").Value = [double]"0.8371556167831109"
$ws.Range("H4").Value = [double]"0.7763661882359298"
$ws.Range("I4").Value = [double]"4.404487885523167E-10"
$ws.Range("J4").Value = [double]"0.000134352010850769"
$ws.Range("B5").Value = [double]"0.2502371920658693"
$ws.Range("C5").Value = [double]"0.003229026081643438"
$ws.Range("D5").Value = [double]"0.02011128779015348"
$ws.Range("F5").Value = [double]"0.2531854706500205"
$ws.Range("G5").Value = [double]"0.018229675389549"
$ws.Range("H5").Value = [double]"0.01035747988784386"
$ws.Range("I5").Value = [double]"0.003408016666508029"
$ws.Range("J5").Value = [double]"0.4017434039354539"
$ws.Range("B6").Value = [double]"0.6387018093091172"
$ws.Range("C6").Value = [double]"0.003849857549612201"
$ws.Range("D6").Value = [double]"0.03618311219722914"
$ws.Range("E6").Value = [double]"0.2531854706500205"
$ws.Range("G6").Value = [double]"0.03556660385828647"
$ws.Range("H6").Value = [double]"0.02014030353133123"
$ws.Range("I6").Value = [double]"0.00404587924091504"
$ws.Range("J6").Value = [double]"0.6216241516277337"
$ws.Range("B7").Value = [double]"0.03847178559399977"
$ws.Range("C7").Value = [double]"7.060412010595485E-07"
$ws.Range("D7").Value = [double]"0.8371556167831109"
$ws.Range("E7").Value = [double]"0.018229675389549"
$ws.Range("F7").Value = [double]"0.03556660385828647"
$ws.Range("H7").Value = [double]"0.8554186857318407"
$ws.Range("I7").Value = [double]"8.772176307569879E-07"
$ws.Range("J7").Value = [double]"6.547176943638711E-05"
$ws.Range("B8").Value = [double]"0.02324467905162675"
$ws.Range("C8").Value = [double]"0.0001782056793002695"
$ws.Range("D8").Value = [double]"0.7763661882359298"
$ws.Range("E8").Value = [double]"0.01035747988784386"
$ws.Range("F8").Value = [double]"0.02014030353133123"
$ws.Range("G8").Value = [double]"0.8554186857318407"
$ws.Range("I8").Value = [double]"0.0001623557350696192"
$ws.Range("J8").Value = [double]"5.687135891774275E-06"
$ws.Range("B9").Value = [double]"0.004331495970726795"
$ws.Range("C9").Value = [double]"0.006419911705607184"
$ws.Range("D9").Value = [double]"4.404487885523167E-10"
$ws.Range("E9").Value = [double]"0.003408016666508029"
$ws.Range("F9").Value = [double]"0.00404587924091504"
$ws.Range("G9").Value = [double]"8.772176307569879E-07"
$ws.Range("H9").Value = [double]"0.0001623557350696192"
$ws.Range("J9").Value = [double]"3.045218133657102E-06"
$ws.Range("B10").Value = [double]"0.504032656816424"
$ws.Range("C10").Value = [double]"4.817164197534751E-06"
$ws.Range("D10").Value = [double]"0.000134352010850769"
$ws.Range("E10").Value = [double]"0.4017434039354539"
$ws.Range("F10").Value = [double]"0.6216241516277337"
$ws.Range("G10").Value = [double]"6.547176943638711E-05"
$ws.Range("H10").Value = [double]"5.687135891774275E-06"
$ws.Range("I10").Value = [double]"3.045218133657102E-06"

$ws = $wb.Worksheets.Item("Estadisticos_DM")
$ws.Range("C2").Value = [double]"3.44554356293951"
$ws.Range("D2").Value = [double]"2.249633324597215"
$ws.Range("E2").Value = [double]"-1.19950975465698"
$ws.Range("F2").Value = [double]"-0.47991131098698"
$ws.Range("G2").Value = [double]"2.284383954478518"
$ws.Range("H2").Value = [double]"2.547168371123652"
$ws.Range("I2").Value = [double]"3.397836052095905"
$ws.Range("J2").Value = [double]"-0.685809349746224"
$ws.Range("B3").Value = [double]"-3.44554356293951"
$ws.Range("D3").Value = [double]"-10.34360809719496"
$ws.Range("E3").Value = [double]"-3.54559908497429"
$ws.Range("F3").Value = [double]"-3.457106733984725"
$ws.Range("G3").Value = [double]"-8.464119253530548"
$ws.Range("H3").Value = [double]"-5.047300102606878"
$ws.Range("I3").Value = [double]"-3.200049298329156"
$ws.Range("J3").Value = [double]"-7.165481192327922"
$ws.Range("B4").Value = [double]"-2.249633324597215"
$ws.Range("C4").Value = [double]"10.34360809719496"
$ws.Range("E4").Value = [double]"-2.621647299533797"
$ws.Range("F4").Value = [double]"-2.316728576027316"
$ws.Range("G4").Value = [double]"-0.2093958397227834"
$ws.Range("H4").Value = [double]"-0.2896009849129749"
$ws.Range("I4").Value = [double]"15.16467443954555"
$ws.Range("J4").Value = [double]"-5.200950035400491"
$ws.Range("B5").Value = [double]"1.19950975465698"
$ws.Range("C5").Value = [double]"3.54559908497429"
$ws.Range("D5").Value = [double]"2.621647299533797"
$ws.Range("F5").Value = [double]"1.191710524936796"
$ws.Range("G5").Value = [double]"2.671966885775691"
$ws.Range("H5").Value = [double]"2.959106641125093"
$ws.Range("I5").Value = [double]"3.518439992896913"
$ws.Range("J5").Value = [double]"0.8647637395590874"
$ws.Range("B6").Value = [double]"0.47991131098698"
$ws.Range("C6").Value = [double]"3.457106733984725"
$ws.Range("D6").Value = [double]"2.316728576027316"
$ws.Range("E6").Value = [double]"-1.191710524936796"
$ws.Range("G6").Value = [double]"2.325771618849315"
$ws.Range("H6").Value = [double]"2.620907572766243"
$ws.Range("I6").Value = [double]"3.432131918402313"
$ws.Range("J6").Value = [double]"-0.5046940976041824"
$ws.Range("B7").Value = [double]"-2.284383954478518"
$ws.Range("C7").Value = [double]"8.464119253530548"
$ws.Range("D7").Value = [double]"0.2093958397227834"
$ws.Range("E7").Value = [double]"-2.671966885775691"
$ws.Range("F7").Value = [double]"-2.325771618849315"
$ws.Range("H7").Value = [double]"-0.185602508983943"
$ws.Range("I7").Value = [double]"8.310067398859305"
$ws.Range("J7").Value = [double]"-5.59993853056273"
$ws.Range("B8").Value = [double]"-2.547168371123652"
$ws.Range("C8").Value = [double]"5.047300102606878"
$ws.Range("D8").Value = [double]"0.2896009849129749"
$ws.Range("E8").Value = [double]"-2.959106641125093"
$ws.Range("F8").Value = [double]"-2.620907572766243"
$ws.Range("G8").Value = [double]"0.185602508983943"
$ws.Range("I8").Value = [double]"5.097783290225642"
$ws.Range("J8").Value = [double]"-7.059576484813658"
$ws.Range("B9").Value = [double]"-3.397836052095905"
$ws.Range("C9").Value = [double]"3.200049298329156"
$ws.Range("D9").Value = [double]"-15.16467443954555"
$ws.Range("E9").Value = [double]"-3.518439992896913"
$ws.Range("F9").Value = [double]"-3.432131918402313"
$ws.Range("G9").Value = [double]"-8.310067398859305"
$ws.Range("H9").Value = [double]"-5.097783290225642"
$ws.Range("J9").Value = [double]"-7.463033488022635"
$ws.Range("B10").Value = [double]"0.685809349746224"
$ws.Range("C10").Value = [double]"7.165481192327922"
$ws.Range("D10").Value = [double]"5.200950035400491"
$ws.Range("E10").Value = [double]"-0.8647637395590874"
$ws.Range("F10").Value = [double]"0.5046940976041824"
$ws.Range("G10").Value = [double]"5.59993853056273"
$ws.Range("H10").Value = [double]"7.059576484813658"
$ws.Range("I10").Value = [double]"7.463033488022635"

$ws = $wb.Worksheets.Item("Resumen")
$ws.Range("B2").Value = [double]"4"
$ws.Range("D2").Value = [double]"4"
$ws.Range("E2").Value = [double]"50"
$ws.Range("B3").Value = [double]"4"
$ws.Range("C3").Value = [double]"0"
$ws.Range("D3").Value = [double]"4"
$ws.Range("E3").Value = [double]"50"
$ws.Range("B4").Value = [double]"1"
$ws.Range("D4").Value = [double]"5"
$ws.Range("E4").Value = [double]"12.5"
$ws.Range("B5").Value = [double]"1"
$ws.Range("D5").Value = [double]"5"
$ws.Range("E5").Value = [double]"12.5"
$ws.Range("B6").Value = [double]"1"
$ws.Range("D6").Value = [double]"5"
$ws.Range("E6").Value = [double]"12.5"
$ws.Range("C7").Value = [double]"0"
$ws.Range("D7").Value = [double]"8"
$ws.Range("C8").Value = [double]"0"
$ws.Range("D8").Value = [double]"8"
$ws.Range("C9").Value = [double]"0"
$ws.Range("D9").Value = [double]"8"

